# Chloride field data workbook - add new sampling-round rows (continuing outlier work)
$wb = $excel.ActiveWorkbook

# --- WIC: new row inserted at row 6 (between 44069 and 44096) ---
$ws = $wb.Worksheets.Item("WIC")
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6, 1).Value = 44082.375
$ws.Cells.Item(6, 2).Value = 56.65
$ws.Cells.Item(6, 3).Value = 17.3

# --- YS: new row inserted at row 19 ---
$ws = $wb.Worksheets.Item("YS")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.489583333336
$ws.Cells.Item(19, 2).Value = 58.99
$ws.Cells.Item(19, 3).Value = 17.2

# --- SW: new row inserted at row 19 ---
$ws = $wb.Worksheets.Item("SW")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.511805555558
$ws.Cells.Item(19, 2).Value = 132.59
$ws.Cells.Item(19, 3).Value = 16.399999999999999

# --- YI: new row inserted at row 19 ---
$ws = $wb.Worksheets.Item("YI")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.39166666667
$ws.Cells.Item(19, 2).Value = 43.71
$ws.Cells.Item(19, 3).Value = 18.600000000000001

# --- YN: new row inserted at row 19 ---
$ws = $wb.Worksheets.Item("YN")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.409722222219
$ws.Cells.Item(19, 2).Value = 35.83
$ws.Cells.Item(19, 3).Value = 17

# --- 6MC: new row inserted at row 19 ---
$ws = $wb.Worksheets.Item("6MC")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.429166666669
$ws.Cells.Item(19, 2).Value = 19.399999999999999
$ws.Cells.Item(19, 3).Value = 14.6

# --- DC: new row inserted at row 19 ---
$ws = $wb.Worksheets.Item("DC")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.439583333333
$ws.Cells.Item(19, 2).Value = 15.66
$ws.Cells.Item(19, 3).Value = 12.9

# --- PBMS: new row inserted at row 19 ---
$ws = $wb.Worksheets.Item("PBMS")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.45416666667
$ws.Cells.Item(19, 2).Value = 81.099999999999994
$ws.Cells.Item(19, 3).Value = 16

# --- PBSF: new row inserted at row 19 (date cell picked up a stray font tweak) ---
$ws = $wb.Worksheets.Item("PBSF")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = 44082.463194444441
$ws.Cells.Item(19, 1).Font.Color = $null
$ws.Cells.Item(19, 2).Value = 51.2
$ws.Cells.Item(19, 3).Value = 16.600000000000001

# --- Update the lingering cell selections left on each sheet after data entry ---
$wb.Worksheets.Item("WIC").Range("C12").Select()
$wb.Worksheets.Item("YS").Range("D34").Select()
$wb.Worksheets.Item("SW").Range("D30").Select()
$wb.Worksheets.Item("YI").Range("B30").Select()
$wb.Worksheets.Item("YN").Range("H20").Select()
$wb.Worksheets.Item("6MC").Range("C27").Select()
$wb.Worksheets.Item("DC").Range("F32").Select()
$wb.Worksheets.Item("PBMS").Range("F27").Select()
$wb.Worksheets.Item("PBSF").Range("F18").Select()
